# Apply the "gh-pages output generated at 456a3b4" update to both the
# "展览" and "全部类型" sheets (they carry identical data tables).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Simple "想去人数" (F column) bumps, same row numbers on both sheets ---
    $ws.Cells.Item(5, 6).Value  = 84     # F5  83  -> 84
    $ws.Cells.Item(6, 6).Value  = 126    # F6  125 -> 126
    $ws.Cells.Item(7, 6).Value  = 1236   # F7  1230 -> 1236
    $ws.Cells.Item(8, 6).Value  = 1522   # F8  1519 -> 1522
    $ws.Cells.Item(9, 6).Value  = 335    # F9  334 -> 335
    $ws.Cells.Item(10, 6).Value = 380    # F10 378 -> 380
    $ws.Cells.Item(16, 6).Value = 271    # F16 270 -> 271
    $ws.Cells.Item(17, 6).Value = 294    # F17 293 -> 294
    $ws.Cells.Item(18, 6).Value = 320    # F18 319 -> 320
    $ws.Cells.Item(19, 6).Value = 1716   # F19 1714 -> 1716
    $ws.Cells.Item(20, 6).Value = 66     # F20 65 -> 66
    $ws.Cells.Item(23, 6).Value = 660    # F23 658 -> 660
    $ws.Cells.Item(26, 6).Value = 4135   # F26 4127 -> 4135
    $ws.Cells.Item(32, 6).Value = 503    # F32 499 -> 503
    $ws.Cells.Item(34, 6).Value = 228    # F34 226 -> 228
    $ws.Cells.Item(36, 6).Value = 134    # F36 133 -> 134

    # --- Rows 28-30: a new event slots in, pushing the old rows 28/29 down
    #     one week, and what was row 30 falls off the bottom (not present
    #     after the edit). Net effect: row28 <- old row29 (new totals),
    #     row29 <- old row30 (new totals), row30 <- old row28 shifted a week
    #     with new totals. ---

    # Row 28 (B28 date stays 2024-02-24)
    $ws.Cells.Item(28, 3).Value = "景德镇·陶溪川×次元文化元宵游园会（ 免费活动）"
    $ws.Cells.Item(28, 4).Value = "新厂西路315号 陶溪川发布大厅"
    $ws.Cells.Item(28, 5).Value = "2024.02.24 10:00-02.25 18:00"
    $ws.Cells.Item(28, 6).Value = 260
    $ws.Cells.Item(28, 7).Value = 30
    $ws.Cells.Item(28, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81207"
    $ws.Cells.Item(28, 9).Value = "//i1.hdslb.com/bfs/openplatform/202402/nIs2jtUn1707298876430.png"

    # Row 29
    # The bare "YYYY-MM-DD" text would otherwise be auto-parsed into a date
    # serial by the COM layer; a leading apostrophe forces text entry, and
    # resetting the style afterwards drops the quote-prefix style so the
    # cell ends up with plain default formatting, same as the source file.
    $ws.Cells.Item(29, 2).Value = "'2024-03-02"
    $ws.Cells.Item(29, 2).Style = "Normal"
    $ws.Cells.Item(29, 3).Value = "南昌·meeting动漫游戏嘉年华"
    $ws.Cells.Item(29, 4).Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
    $ws.Cells.Item(29, 5).Value = "2024.03.02 09:00-03.03 17:00"
    $ws.Cells.Item(29, 6).Value = 1077
    $ws.Cells.Item(29, 7).Value = 60
    $ws.Cells.Item(29, 8).Value = "https://show.bilibili.com/platform/detail.html?id=79555"
    $ws.Cells.Item(29, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/l6GUtggC1706843695971.jpeg"

    # Row 30
    $ws.Cells.Item(30, 2).Value = "'2024-03-09"
    $ws.Cells.Item(30, 2).Style = "Normal"
    $ws.Cells.Item(30, 3).Value = "景德镇·江报国风动漫展 "
    $ws.Cells.Item(30, 4).Value = "迎宾大道与寺山路交叉口东200米 陶博城"
    $ws.Cells.Item(30, 5).Value = "2024.03.09 09:00-03.10 17:00"
    $ws.Cells.Item(30, 6).Value = 479
    $ws.Cells.Item(30, 7).Value = 45
    $ws.Cells.Item(30, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81362"
    $ws.Cells.Item(30, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/ae5G3ouV1706092057911.jpeg"
}
